# chore: update Sheets via scheduled runner
# Refresh market-price-derived figures (currentAveragePrice*, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) across the per-job Leve profit sheets with newly pulled
# market board data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 129 (ALC)
$ws.Range("H129").Value = 1041.098
$ws.Range("J129").Value = 1131.5682
$ws.Range("L129").Value = 3394.7046
$ws.Range("N129").Value = -13394.7046

# Row 132 (ALC)
$ws.Range("H132").Value = 1829.2972
$ws.Range("I132").Value = 1946.1515
$ws.Range("J132").Value = 865.25
$ws.Range("K132").Value = 5838.4545
$ws.Range("L132").Value = 2595.75
$ws.Range("M132").Value = -3308.4545
$ws.Range("N132").Value = -7655.75

# Row 135 (ALC)
$ws.Range("H135").Value = 670.7406999999999
$ws.Range("I135").Value = 670.7406999999999
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 6036.6663
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3501.6663
$ws.Range("N135").ClearContents()

# Row 137 (ALC)
$ws.Range("H137").Value = 1409.9767
$ws.Range("I137").Value = 1340.3928
$ws.Range("J137").Value = 1539.8667
$ws.Range("K137").Value = 4021.1784
$ws.Range("L137").Value = 4619.6001
$ws.Range("M137").Value = -1471.1784
$ws.Range("N137").Value = -9719.6001

# Row 138 (ALC)
$ws.Range("H138").Value = 3419
$ws.Range("I138").Value = 1747.2667
$ws.Range("J138").Value = 4672.8
$ws.Range("K138").Value = 5241.800099999999
$ws.Range("L138").Value = 14018.4
$ws.Range("M138").Value = -101.8000999999995
$ws.Range("N138").Value = -24298.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (ARM)
$ws.Range("H32").Value = 3522.2559
$ws.Range("I32").Value = 2286.7297
$ws.Range("J32").Value = 11141.333
$ws.Range("K32").Value = 2286.7297
$ws.Range("L32").Value = 11141.333
$ws.Range("M32").Value = -1999.7297
$ws.Range("N32").Value = -11715.333

# Row 74 (ARM)
$ws.Range("H74").Value = 1091.0178
$ws.Range("I74").Value = 1089.1957
$ws.Range("J74").Value = 1099.4
$ws.Range("K74").Value = 1089.1957
$ws.Range("L74").Value = 1099.4
$ws.Range("M74").Value = -215.1957
$ws.Range("N74").Value = -2847.4

# Row 77 (ARM)
$ws.Range("H77").Value = 1091.0178
$ws.Range("I77").Value = 1089.1957
$ws.Range("J77").Value = 1099.4
$ws.Range("K77").Value = 5445.9785
$ws.Range("L77").Value = 5497
$ws.Range("M77").Value = -1077.9785
$ws.Range("N77").Value = -14233

$ws = $wb.Worksheets.Item("BSM")
# Row 134 (BSM)
$ws.Range("H134").Value = 1722.3462
$ws.Range("I134").Value = 808.7895
$ws.Range("K134").Value = 2426.3685
$ws.Range("M134").Value = 108.6315

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 3874.1643
$ws.Range("I31").Value = 789.8158
$ws.Range("J31").Value = 7222.8857
$ws.Range("K31").Value = 789.8158
$ws.Range("L31").Value = 7222.8857
$ws.Range("M31").Value = -494.8158
$ws.Range("N31").Value = -7812.8857

# Row 34 (CRP)
$ws.Range("H34").Value = 3874.1643
$ws.Range("I34").Value = 789.8158
$ws.Range("J34").Value = 7222.8857
$ws.Range("K34").Value = 789.8158
$ws.Range("L34").Value = 7222.8857
$ws.Range("M34").Value = -587.8158
$ws.Range("N34").Value = -7626.8857

# Row 134 (CRP)
$ws.Range("H134").Value = 4554.4546
$ws.Range("I134").Value = 4851.6553
$ws.Range("K134").Value = 14554.9659
$ws.Range("M134").Value = -12019.9659

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 283802.75
$ws.Range("I5").Value = 359.56
$ws.Range("J5").Value = 475318.44
$ws.Range("K5").Value = 1078.68
$ws.Range("L5").Value = 1425955.32
$ws.Range("M5").Value = -966.6800000000001
$ws.Range("N5").Value = -1426179.32

# Row 68 (CUL)
$ws.Range("H68").Value = 564057.5
$ws.Range("I68").Value = 1666519
$ws.Range("J68").Value = 1577.1632
$ws.Range("K68").Value = 4999557
$ws.Range("L68").Value = 4731.4896
$ws.Range("M68").Value = -4998746
$ws.Range("N68").Value = -6353.4896

# Row 71 (CUL)
$ws.Range("H71").Value = 564057.5
$ws.Range("I71").Value = 1666519
$ws.Range("J71").Value = 1577.1632
$ws.Range("K71").Value = 14998671
$ws.Range("L71").Value = 14194.4688
$ws.Range("M71").Value = -14994615
$ws.Range("N71").Value = -22306.4688

# Row 113 (CUL)
$ws.Range("H113").Value = 1782.48
$ws.Range("I113").Value = 2828.75
$ws.Range("J113").Value = 1290.1177
$ws.Range("K113").Value = 8486.25
$ws.Range("L113").Value = 3870.3531
$ws.Range("M113").Value = -6316.25
$ws.Range("N113").Value = -8210.3531

# Row 131 (CUL)
$ws.Range("H131").Value = 4894.241
$ws.Range("I131").Value = 1050
$ws.Range("J131").Value = 5179
$ws.Range("K131").Value = 3150
$ws.Range("L131").Value = 15537
$ws.Range("M131").Value = 1890
$ws.Range("N131").Value = -25617

# Row 132 (CUL)
$ws.Range("H132").Value = 1650
$ws.Range("I132").Value = 1757.6923
$ws.Range("J132").Value = 1494.4445
$ws.Range("K132").Value = 15819.2307
$ws.Range("L132").Value = 13450.0005
$ws.Range("M132").Value = -13289.2307
$ws.Range("N132").Value = -18510.0005

# Row 135 (CUL)
$ws.Range("H135").Value = 283802.75
$ws.Range("I135").Value = 359.56
$ws.Range("J135").Value = 475318.44
$ws.Range("K135").Value = 3236.04
$ws.Range("L135").Value = 4277865.96
$ws.Range("M135").Value = -701.04
$ws.Range("N135").Value = -4282935.96

$ws = $wb.Worksheets.Item("GSM")
# Row 63 (GSM)
$ws.Range("H63").Value = 6500
$ws.Range("J63").Value = 6500
$ws.Range("L63").Value = 6500
$ws.Range("N63").Value = -7872

# Row 66 (GSM)
$ws.Range("H66").Value = 6500
$ws.Range("J66").Value = 6500
$ws.Range("L66").Value = 19500
$ws.Range("N66").Value = -26364

# Row 68 (GSM)
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

# Row 71 (GSM)
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# Row 86 (GSM)
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89 (GSM)
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 132 (GSM)
$ws.Range("H132").Value = 6900.75
$ws.Range("I132").Value = 5640
$ws.Range("K132").Value = 16920
$ws.Range("M132").Value = -14390

$ws = $wb.Worksheets.Item("LTW")
# Row 132 (LTW)
$ws.Range("H132").Value = 5745.0625
$ws.Range("I132").Value = 8036.231
$ws.Range("J132").Value = 4177.421
$ws.Range("K132").Value = 24108.693
$ws.Range("L132").Value = 12532.263
$ws.Range("M132").Value = -21578.693
$ws.Range("N132").Value = -17592.263

# Row 136 (LTW)
$ws.Range("H136").Value = 16670741
$ws.Range("I136").Value = 3775.3333
$ws.Range("K136").Value = 11325.9999
$ws.Range("M136").Value = -8775.999899999999
